$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missing_years")

# Row 22 (Ohio University) gained Contacted?/Notes/Completed? = 1 values
$ws.Range("B22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1

# New rows 46-50: five universities appended to the missing_years log
$ws.Range("A46").Value = "Texas Christian University"
$ws.Range("B46").Value = 1
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0

$ws.Range("A47").Value = "Hampden-Sydney College"
$ws.Range("B47").Value = 1
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
# Row 47's university name keeps the same highlighted style as rows 2-21
$ws.Range("A21").Copy()
$ws.Range("A47").PasteSpecial(-4122)

$ws.Range("A48").Value = "University of Kentucky"
$ws.Range("B48").Value = 1
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0

$ws.Range("A49").Value = "University of Nebraska-Lincoln"
$ws.Range("B49").Value = 1
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 1
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0

$ws.Range("A50").Value = "University of Nevada-Reno"
$ws.Range("B50").Value = 1
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 1
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0

# Move the active tab / selection from Sheet1 to missing_years
$ws.Activate()
$ws.Range("B45").Select()
